$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$tl = $s.TimeLine
$seq = $tl.MainSequence
Write-Output $seq.Count
